# V 2.0.2 se arreglo la fechar y hora de reimpresion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Patient identification
$ws.Range("A6").Value = "CACHUC"
$ws.Range("C6").Value = "LOPEZ"
$ws.Range("E6").Value = "MARIO "
$ws.Range("G6").Value = "FELIPE "
$ws.Range("I6").Value = "4421/201755006"

# Date of birth / age / place of birth / sex
$ws.Range("A12").Value = "1977-05-05"
$ws.Range("F12").Value = "40"
$ws.Range("H12").Value = "GUATEMALA "
$ws.Range("J12").Value = "null"

# Estado civil / ocupacion / nacionalidad / no. cedula
$ws.Range("A14").Value = "null"
$ws.Range("D14").Value = ""
$ws.Range("H14").Value = ""

# Emergency contact
$ws.Range("A20").Value = "LUCIA LOPEZ "
$ws.Range("F20").Value = "MAMA "
$ws.Range("H20").Value = "LOTE 23 COL. 10 MAYO Z. 7"
$ws.Range("J20").Value = "50835942"

# Reprint date / time fix
$ws.Range("A24").Value = "24/10/2017"
$ws.Range("C24").Value = "15:9:49"
